$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vegfc"
$ws.Range("C2").Value = "Nrp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 3.911257666666666
$ws.Range("H2").Value = 11.733773
$ws.Range("I2").Value = 0.4115343446855154
$ws.Range("J2").Value = 0.4115343446855154
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 34.52052333333333
$ws.Range("N2").Value = 103.56157
$ws.Range("O2").Value = 0.7684334662422598
$ws.Range("P2").Value = 0.7684334662422598
$ws.Range("Q2").Value = 135.0186615448455
$ws.Range("R2").Value = 1215.16795390361
$ws.Range("S2").Value = 0.3162367629644275
$ws.Range("T2").Value = 0.3162367629644275

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vegfc"
$ws.Range("C3").Value = "Nrp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 3.911257666666666
$ws.Range("H3").Value = 11.733773
$ws.Range("I3").Value = 0.4115343446855154
$ws.Range("J3").Value = 0.4115343446855154
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.347618333333333
$ws.Range("N3").Value = 16.042855
$ws.Range("O3").Value = 0.1190390091234806
$ws.Range("P3").Value = 0.1190390091234805
$ws.Range("Q3").Value = 20.91591320465722
$ws.Range("R3").Value = 188.243218841915
$ws.Range("S3").Value = 0.04898864061164466
$ws.Range("T3").Value = 0.04898864061164465

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vegfc"
$ws.Range("C4").Value = "Nrp2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 3.911257666666666
$ws.Range("H4").Value = 11.733773
$ws.Range("I4").Value = 0.4115343446855154
$ws.Range("J4").Value = 0.4115343446855154
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.055101333333334
$ws.Range("N4").Value = 15.165304
$ws.Range("O4").Value = 0.1125275246342597
$ws.Range("P4").Value = 0.1125275246342597
$ws.Range("Q4").Value = 19.77180384577689
$ws.Range("R4").Value = 177.946234611992
$ws.Range("S4").Value = 0.04630894110944325
$ws.Range("T4").Value = 0.04630894110944325

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vegfc"
$ws.Range("C5").Value = "Nrp2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.100181333333333
$ws.Range("H5").Value = 12.300544
$ws.Range("I5").Value = 0.4314124974392592
$ws.Range("J5").Value = 0.4314124974392592
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 34.52052333333333
$ws.Range("N5").Value = 103.56157
$ws.Range("O5").Value = 0.7684334662422598
$ws.Range("P5").Value = 0.7684334662422598
$ws.Range("Q5").Value = 141.5404053882311
$ws.Range("R5").Value = 1273.86364849408
$ws.Range("S5").Value = 0.33151180078748
$ws.Range("T5").Value = 0.33151180078748

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vegfc"
$ws.Range("C6").Value = "Nrp2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.100181333333333
$ws.Range("H6").Value = 12.300544
$ws.Range("I6").Value = 0.4314124974392592
$ws.Range("J6").Value = 0.4314124974392592
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.347618333333333
$ws.Range("N6").Value = 16.042855
$ws.Range("O6").Value = 0.1190390091234806
$ws.Range("P6").Value = 0.1190390091234805
$ws.Range("Q6").Value = 21.92620486812444
$ws.Range("R6").Value = 197.33584381312
$ws.Range("S6").Value = 0.05135491621865551
$ws.Range("T6").Value = 0.0513549162186555

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vegfc"
$ws.Range("C7").Value = "Nrp2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.100181333333333
$ws.Range("H7").Value = 12.300544
$ws.Range("I7").Value = 0.4314124974392592
$ws.Range("J7").Value = 0.4314124974392592
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.055101333333334
$ws.Range("N7").Value = 15.165304
$ws.Range("O7").Value = 0.1125275246342597
$ws.Range("P7").Value = 0.1125275246342597
$ws.Range("Q7").Value = 20.72683212504178
$ws.Range("R7").Value = 186.541489125376
$ws.Range("S7").Value = 0.04854578043312374
$ws.Range("T7").Value = 0.04854578043312373

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Vegfc"
$ws.Range("C8").Value = "Nrp2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.492646666666667
$ws.Range("H8").Value = 4.47794
$ws.Range("I8").Value = 0.1570531578752254
$ws.Range("J8").Value = 0.1570531578752254
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 34.52052333333333
$ws.Range("N8").Value = 103.56157
$ws.Range("O8").Value = 0.7684334662422598
$ws.Range("P8").Value = 0.7684334662422598
$ws.Range("Q8").Value = 51.52694408508889
$ws.Range("R8").Value = 463.7424967658
$ws.Range("S8").Value = 0.1206849024903523
$ws.Range("T8").Value = 0.1206849024903523

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Vegfc"
$ws.Range("C9").Value = "Nrp2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.492646666666667
$ws.Range("H9").Value = 4.47794
$ws.Range("I9").Value = 0.1570531578752254
$ws.Range("J9").Value = 0.1570531578752254
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.347618333333333
$ws.Range("N9").Value = 16.042855
$ws.Range("O9").Value = 0.1190390091234806
$ws.Range("P9").Value = 0.1190390091234805
$ws.Range("Q9").Value = 7.982104679855555
$ws.Range("R9").Value = 71.8389421187
$ws.Range("S9").Value = 0.01869545229318039
$ws.Range("T9").Value = 0.01869545229318039

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vegfc"
$ws.Range("C10").Value = "Nrp2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.492646666666667
$ws.Range("H10").Value = 4.47794
$ws.Range("I10").Value = 0.1570531578752254
$ws.Range("J10").Value = 0.1570531578752254
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.055101333333334
$ws.Range("N10").Value = 15.165304
$ws.Range("O10").Value = 0.1125275246342597
$ws.Range("P10").Value = 0.1125275246342597
$ws.Range("Q10").Value = 7.545480154862223
$ws.Range("R10").Value = 67.90932139376001
$ws.Range("S10").Value = 0.0176728030916927
$ws.Range("T10").Value = 0.0176728030916927
